$wb = $excel.ActiveWorkbook

# --- Sheet "summary" ---
$ws1 = $wb.Worksheets.Item("summary")

$ws1.Range("C2").Value = 680
$ws1.Range("D2").Value = 644
$ws1.Range("E2").Value = 77.02
$ws1.Range("F2").Value = -1.44
$ws1.Range("I2").Value = -0.07
$ws1.Range("L2").Value = 1.07

$ws1.Range("C3").Value = 680
$ws1.Range("D3").Value = 575
$ws1.Range("E3").Value = 41.74
$ws1.Range("F3").Value = 0.42
$ws1.Range("H3").Value = 0.99
$ws1.Range("I3").Value = -0.22
$ws1.Range("J3").Value = 0.29
$ws1.Range("L3").Value = 1.08

$ws1.Range("C4").Value = 680
$ws1.Range("D4").Value = 349
$ws1.Range("E4").Value = 12.03
$ws1.Range("F4").Value = 2.36
$ws1.Range("G4").Value = 0.17
$ws1.Range("H4").Value = 0.98
$ws1.Range("I4").Value = -0.14
$ws1.Range("J4").Value = 0.24
$ws1.Range("L4").Value = 1.5

$ws1.Range("C5").Value = 680
$ws1.Range("D5").Value = 657
$ws1.Range("E5").Value = 75.95
$ws1.Range("F5").Value = -1.36
$ws1.Range("H5").Value = 1.01
$ws1.Range("I5").Value = 0.18
$ws1.Range("J5").Value = 0.22
$ws1.Range("K5").Value = 0.03
$ws1.Range("L5").Value = 0.89

$ws1.Range("C6").Value = 680
$ws1.Range("D6").Value = 655
$ws1.Range("E6").Value = 72.52
$ws1.Range("G6").Value = 0.09
$ws1.Range("L6").Value = 0.9

$ws1.Range("C7").Value = 680
$ws1.Range("D7").Value = 647
$ws1.Range("E7").Value = 64.91
$ws1.Range("F7").Value = -0.74
$ws1.Range("H7").Value = 0.99
$ws1.Range("I7").Value = -0.25
$ws1.Range("J7").Value = 0.28
$ws1.Range("L7").Value = 1.06

$ws1.Range("C8").Value = 680
$ws1.Range("D8").Value = 648
$ws1.Range("E8").Value = 62.35
$ws1.Range("F8").Value = -0.61
$ws1.Range("H8").Value = 1.02
$ws1.Range("I8").Value = 0.45
$ws1.Range("K8").Value = 0.03
$ws1.Range("L8").Value = 0.83

# --- Sheet "model_fit" ---
$ws2 = $wb.Worksheets.Item("model_fit")

$ws2.Range("B2").Value = 680
$ws2.Range("D2").Value = 4761
$ws2.Range("E2").Value = 4777
$ws2.Range("F2").Value = 4813
$ws2.Range("G2").Value = 0.497
$ws2.Range("H2").Value = 0.281

$ws2.Range("B3").Value = 680
$ws2.Range("D3").Value = 4757
$ws2.Range("E3").Value = 4785
$ws2.Range("F3").Value = 4849
$ws2.Range("G3").Value = 0.505
$ws2.Range("H3").Value = 0.293
